$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eutrophication")

# Copy the (already-formatted) row 5 cells' format onto rows 6-8 so all
# four rows share the same style (matches the xf index consolidation seen
# in the diff: C6:E8 go from s="6" to s="4", same as C5:E5).
$ws.Range("C5:E5").Copy()
$ws.Range("C6:E8").PasteSpecial(-4122)

# Fill in the previously-blank cells with 0 (LCA design completed).
$ws.Range("C5:E8").Value = 0

# Make "Eutrophication" the active sheet/tab (was "info" before) and move
# the selection on this sheet to D13.
$ws.Activate() | Out-Null
$ws.Range("D13").Select() | Out-Null
